$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Variable")
$ws2 = $wb.Worksheets.Item("r Variable_DataType")

$ws1.Rows.Item(30).Delete()
$ws2.Rows.Item(30).Delete()
